$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1978.5555
$ws.Range("I15").Value = 1978.5555
$ws.Range("K15").Value = 5935.666499999999
$ws.Range("M15").Value = -5766.666499999999

# Row 33
$ws.Range("H33").Value = 219.44444

# Row 96
$ws.Range("H96").Value = 8349.833000000001
$ws.Range("J96").Value = 2019.8
$ws.Range("L96").Value = 6059.4
$ws.Range("N96").Value = -8805.4

# Row 100
$ws.Range("H100").Value = 838.5
$ws.Range("I100").Value = 323.625
$ws.Range("K100").Value = 323.625
$ws.Range("M100").Value = 217.375

# Row 106
$ws.Range("H106").Value = 2006
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

# Row 107
$ws.Range("H107").Value = 278.7143
$ws.Range("I107").Value = 169.8
$ws.Range("J107").Value = 551
$ws.Range("K107").Value = 169.8
$ws.Range("L107").Value = 551
$ws.Range("M107").Value = 1750.2
$ws.Range("N107").Value = -4391

# Row 116
$ws.Range("H116").Value = 3661.72
$ws.Range("I116").Value = 3059
$ws.Range("K116").Value = 3059
$ws.Range("M116").Value = 383

# Row 138
$ws.Range("H138").Value = 4852.077
$ws.Range("I138").Value = 3763.1667
$ws.Range("J138").Value = 5785.4287
$ws.Range("K138").Value = 11289.5001
$ws.Range("L138").Value = 17356.2861
$ws.Range("M138").Value = -6149.500100000001
$ws.Range("N138").Value = -27636.2861

$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 1763.4117
$ws.Range("I110").Value = 1763.4117
$ws.Range("K110").Value = 1763.4117
$ws.Range("M110").Value = 281.5882999999999

# Row 132
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 5234.3335
$ws.Range("I20").Value = 4760.636
$ws.Range("K20").Value = 4760.636
$ws.Range("M20").Value = -4513.636

# Row 99
$ws.Range("H99").Value = 2902
$ws.Range("I99").Value = 2836.6667
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2836.6667
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -1338.6667
$ws.Range("N99").Value = -5996

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1791.125
$ws.Range("J31").Value = 2099.3333
$ws.Range("L31").Value = 2099.3333
$ws.Range("N31").Value = -2689.3333

# Row 34
$ws.Range("H34").Value = 1791.125
$ws.Range("J34").Value = 2099.3333
$ws.Range("L34").Value = 2099.3333
$ws.Range("N34").Value = -2503.3333

# Row 74
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# Row 86
$ws.Range("H86").Value = 4949.25
$ws.Range("I86").Value = 4950
$ws.Range("J86").Value = 4949
$ws.Range("K86").Value = 4950
$ws.Range("L86").Value = 4949
$ws.Range("M86").Value = -3827
$ws.Range("N86").Value = -7195

# Row 88
$ws.Range("H88").Value = 19177.572
$ws.Range("J88").Value = 19177.572
$ws.Range("L88").Value = 19177.572
$ws.Range("N88").Value = -19989.572

# Row 89
$ws.Range("H89").Value = 4949.25
$ws.Range("I89").Value = 4950
$ws.Range("J89").Value = 4949
$ws.Range("K89").Value = 24750
$ws.Range("L89").Value = 24745
$ws.Range("M89").Value = -19134
$ws.Range("N89").Value = -35977

# Row 91
$ws.Range("H91").Value = 19177.572
$ws.Range("J91").Value = 19177.572
$ws.Range("L91").Value = 19177.572
$ws.Range("N91").Value = -21985.572

# Row 94
$ws.Range("H94").Value = 126515.336
$ws.Range("I94").Value = 278477.75
$ws.Range("K94").Value = 278477.75
$ws.Range("M94").Value = -278026.75

# Row 96
$ws.Range("H96").Value = 29299.75
$ws.Range("J96").Value = 29299.75
$ws.Range("L96").Value = 29299.75
$ws.Range("N96").Value = -34791.75

# Row 99
$ws.Range("H99").Value = 1200
$ws.Range("I99").Value = 1200
$ws.Range("K99").Value = 1200
$ws.Range("M99").Value = 298

# Row 122
$ws.Range("H122").Value = 5999
$ws.Range("I122").Value = 1999
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 5997
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -3547
$ws.Range("N122").Value = -34897

# Row 126
$ws.Range("H126").Value = 1200
$ws.Range("I126").Value = 1200
$ws.Range("K126").Value = 3600
$ws.Range("M126").Value = -1130

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 3972.5833
$ws.Range("I7").Value = 4332.8184
$ws.Range("J7").Value = 10
$ws.Range("K7").Value = 12998.4552
$ws.Range("L7").Value = 30
$ws.Range("M7").Value = -12886.4552
$ws.Range("N7").Value = -254

# Row 76
$ws.Range("H76").Value = 12201.2
$ws.Range("I76").Value = 6202.4
$ws.Range("K76").Value = 18607.2
$ws.Range("M76").Value = -18224.2

# Row 79
$ws.Range("H79").Value = 12201.2
$ws.Range("I79").Value = 6202.4
$ws.Range("K79").Value = 18607.2
$ws.Range("M79").Value = -17281.2

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 126.545456
$ws.Range("I2").Value = 133.05556
$ws.Range("K2").Value = 133.05556
$ws.Range("M2").Value = -20.05556000000001

# Row 35
$ws.Range("H35").Value = 22500
$ws.Range("J35").Value = 10000
$ws.Range("L35").Value = 10000
$ws.Range("N35").Value = -10596

# Row 57
$ws.Range("H57").Value = 78579.39999999999
$ws.Range("J57").Value = 78579.39999999999
$ws.Range("L57").Value = 78579.39999999999
$ws.Range("N57").Value = -80219.39999999999

# Row 102
$ws.Range("H102").Value = 1863.8334
$ws.Range("I102").Value = 1837
$ws.Range("K102").Value = 1837
$ws.Range("M102").Value = -215

# Row 122
$ws.Range("H122").Value = 2533.5
$ws.Range("I122").Value = 2424.7273
$ws.Range("K122").Value = 7274.1819
$ws.Range("M122").Value = -4824.1819

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1775.1666
$ws.Range("J22").Value = 2922
$ws.Range("L22").Value = 2922
$ws.Range("N22").Value = -3512

# Row 27
$ws.Range("H27").Value = 1775.1666
$ws.Range("J27").Value = 2922
$ws.Range("L27").Value = 2922
$ws.Range("N27").Value = -3136

# Row 40
$ws.Range("H40").Value = 2910.7693
$ws.Range("I40").Value = 2078.4285
$ws.Range("K40").Value = 2078.4285
$ws.Range("M40").Value = -1942.4285

# Row 55
$ws.Range("H55").Value = 1316.8422
$ws.Range("I55").Value = 1139.5454
$ws.Range("K55").Value = 1139.5454
$ws.Range("M55").Value = -966.5454

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 4603.727
$ws.Range("I122").Value = 2724.4285
$ws.Range("K122").Value = 8173.2855
$ws.Range("M122").Value = -5723.2855

# Row 126
$ws.Range("H126").Value = 1406.2858
$ws.Range("I126").Value = 1406.2858
$ws.Range("K126").Value = 4218.857400000001
$ws.Range("M126").Value = -1748.857400000001
